$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row with the new sample value that contains a comma (quoted
# literally, matching the existing data's pattern of leading/escaping quotes).
$ws.Range("A6").Value = '"parent,06"'

# Update the active selection as recorded in the saved workbook.
$ws.Range("C8").Select()
